$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Jane"
$ws.Range("C2").Value = "Doe"
$ws.Range("D2").Value = "demo@email.com"
$ws.Range("E2").Value = 109.4
$ws.Range("F2").Value = 23
$ws.Range("G2").Value = "Loss of Movement,Fever,Tiredness"
$ws.Range("H2").Value = 3
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 2
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = $false
$ws.Range("M2").Value = "Very High Risk"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Jake"
$ws.Range("C3").Value = "Doe"
$ws.Range("D3").Value = "demo@email.com"
$ws.Range("E3").Value = 73.40000000000001
$ws.Range("F3").Value = 34
$ws.Range("G3").Value = "Difficulty Breathing,Tiredness"
$ws.Range("H3").Value = 3
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = $false
$ws.Range("M3").Value = "Not at Risk"
